# Auto-generated PowerPoint COM-interop edit script
# Rewrites the 'Education' Wikipedia-scrape deck into a 'Wind' deck,
# and appends 20 new slides (40-59) pulled from the Wind Wikipedia article/menus.

$p = $ppt.ActivePresentation

# --- Slide 2 (title slide) ---
$p.Slides.Item(2).Shapes.Item(1).TextFrame.TextRange.Text = "Wind"
$p.Slides.Item(2).Shapes.Item(2).TextFrame.TextRange.Text = " by aman"

# --- Slides 3-6 (table-of-contents bullet lists) ---
$p.Slides.Item(3).Shapes.Item(2).TextFrame.TextRange.Text = " Contents`r Causes[edit]`r Measurement[edit]`r Wind force scale[edit]`r  -- Enhanced Fujita scale[edit]`r  -- Station model[edit]`r Wind power[edit]`r  -- Theoretical power captured by a wind turbine[edit]`r  -- Practical wind turbine power[edit]`r"
$p.Slides.Item(4).Shapes.Item(2).TextFrame.TextRange.Text = " Global climatology[edit]`r  -- Tropics[edit]`r  -- Westerlies and their impact[edit]`r  -- Polar easterlies[edit]`r Local considerations[edit]`r  -- Sea and land breezes[edit]`r  -- Near mountains[edit]`r Average wind speeds[edit]`r"
$p.Slides.Item(5).Shapes.Item(2).TextFrame.TextRange.Text = " Shear[edit]`r Usage[edit]`r  -- History[edit]`r  -- Transportation[edit]`r  -- Power source[edit]`r  -- Recreation[edit]`r Role in the natural world[edit]`r  -- Erosion[edit]`r  -- Desert dust migration[edit]`r"
$p.Slides.Item(6).Shapes.Item(2).TextFrame.TextRange.Text = "  -- Effect on plants[edit]`r  -- Effect on animals[edit]`r  -- Sound generation[edit]`r Related damage[edit]`r In outer space[edit]`r  -- Planetary wind[edit]`r  -- Solar wind[edit]`r On other planets[edit]`r See also[edit]`r"

# --- Slides 7-9 (leftover scratch TOC / first content slides) ---
$p.Slides.Item(7).Shapes.Item(1).TextFrame.TextRange.Text = "contents"
$p.Slides.Item(7).Shapes.Item(2).TextFrame.TextRange.Text = " References[edit]`r External links[edit]`r Navigation menu`r  -- Personal tools`r  -- Namespaces`r  -- Variants`r  -- Views`r  -- More`r  -- Search`r"
$p.Slides.Item(8).Shapes.Item(1).TextFrame.TextRange.Text = "contents"
$p.Slides.Item(8).Shapes.Item(2).TextFrame.TextRange.Text = "  -- Navigation`r  -- Contribute`r  -- Tools`r  -- Print/export`r  -- In other projects`r  -- Languages`r"
$p.Slides.Item(9).Shapes.Item(1).TextFrame.TextRange.Text = "Contents"
$p.Slides.Item(9).Shapes.Item(2).TextFrame.TextRange.Text = "Contentsh2 Causes[edit]p The wind is caused by differences in atmospheric pressure."

# --- Slides 10-39 (content slides: Education topics -> Wind topics) ---
$p.Slides.Item(10).Shapes.Item(1).TextFrame.TextRange.Text = "Causes[edit]"
$p.Slides.Item(10).Shapes.Item(2).TextFrame.TextRange.Text = "Causes[edit]p The wind is caused by differences in atmospheric pressure."
$p.Slides.Item(11).Shapes.Item(1).TextFrame.TextRange.Text = "Measurement[edit]"
$p.Slides.Item(11).Shapes.Item(2).TextFrame.TextRange.Text = "Measurement[edit]p Wind direction is usually expressed in terms of the direction from which it originates."
$p.Slides.Item(12).Shapes.Item(1).TextFrame.TextRange.Text = "Wind force scale[edit]"
$p.Slides.Item(12).Shapes.Item(2).TextFrame.TextRange.Text = "Wind force scale[edit]p Historically, the Beaufort wind force scale (created by Beaufort) provides an empirical description of wind speed based on observed sea conditions."
$p.Slides.Item(13).Shapes.Item(1).TextFrame.TextRange.Text = "Enhanced Fujita scale[edit]"
$p.Slides.Item(13).Shapes.Item(2).TextFrame.TextRange.Text = "Enhanced Fujita scale[edit]p The Enhanced Fujita Scale (EF Scale) rates the strength of tornadoes in the United States by using damage to estimate wind speed."
$p.Slides.Item(14).Shapes.Item(1).TextFrame.TextRange.Text = "Station model[edit]"
$p.Slides.Item(14).Shapes.Item(2).TextFrame.TextRange.Text = "Station model[edit]p The station model plotted on surface weather maps uses a wind barb to show both wind direction and speed."
$p.Slides.Item(15).Shapes.Item(1).TextFrame.TextRange.Text = "Wind power[edit]"
$p.Slides.Item(15).Shapes.Item(2).TextFrame.TextRange.Text = "Wind power[edit]p Wind energy is the kinetic energy of the air in motion."
$p.Slides.Item(16).Shapes.Item(1).TextFrame.TextRange.Text = "Theoretical power captured by a wind turbine[edit]"
$p.Slides.Item(16).Shapes.Item(2).TextFrame.TextRange.Text = "Theoretical power captured by a wind turbine[edit]p Total wind power could be captured only if the wind velocity is reduced to zero."
$p.Slides.Item(17).Shapes.Item(1).TextFrame.TextRange.Text = "Practical wind turbine power[edit]"
$p.Slides.Item(17).Shapes.Item(2).TextFrame.TextRange.Text = "Practical wind turbine power[edit]p Further insufficiencies, such as rotor blade friction and drag, gearbox losses, generator and converter losses, reduce the power delivered by a wind turbine."
$p.Slides.Item(18).Shapes.Item(1).TextFrame.TextRange.Text = "Global climatology[edit]"
$p.Slides.Item(18).Shapes.Item(2).TextFrame.TextRange.Text = "Global climatology[edit]p Easterly winds, on average, dominate the flow pattern across the poles, westerly winds blow across the mid-latitudes of the earth, polewards of the subtropical ridge, while easterlies again dominate the tropics."
$p.Slides.Item(19).Shapes.Item(1).TextFrame.TextRange.Text = "Tropics[edit]"
$p.Slides.Item(19).Shapes.Item(2).TextFrame.TextRange.Text = "Tropics[edit]p The trade winds (also called trades) are the prevailing pattern of easterly surface winds found in the tropics towards the Earth's equator."
$p.Slides.Item(20).Shapes.Item(1).TextFrame.TextRange.Text = "Westerlies and their impact[edit]"
$p.Slides.Item(20).Shapes.Item(2).TextFrame.TextRange.Text = "Westerlies and their impact[edit]p The Westerlies or the Prevailing Westerlies are the prevailing winds in the middle latitudes between 35 and 65 degrees latitude."
$p.Slides.Item(21).Shapes.Item(1).TextFrame.TextRange.Text = "Polar easterlies[edit]"
$p.Slides.Item(21).Shapes.Item(2).TextFrame.TextRange.Text = "Polar easterlies[edit]p The polar easterlies, also known as Polar Hadley cells, are dry, cold prevailing winds that blow from the high-pressure areas of the polar highs at the north and south poles towards the low-pressure areas within the Westerlies at high latitudes."
$p.Slides.Item(22).Shapes.Item(1).TextFrame.TextRange.Text = "Local considerations[edit]"
$p.Slides.Item(22).Shapes.Item(2).TextFrame.TextRange.Text = "Local considerations[edit]h3 Sea and land breezes[edit]p In coastal regions, sea breezes and land breezes can be important factors in a location's prevailing winds."
$p.Slides.Item(23).Shapes.Item(1).TextFrame.TextRange.Text = "Sea and land breezes[edit]"
$p.Slides.Item(23).Shapes.Item(2).TextFrame.TextRange.Text = "Sea and land breezes[edit]p In coastal regions, sea breezes and land breezes can be important factors in a location's prevailing winds."
$p.Slides.Item(24).Shapes.Item(1).TextFrame.TextRange.Text = "Near mountains[edit]"
$p.Slides.Item(24).Shapes.Item(2).TextFrame.TextRange.Text = "Near mountains[edit]p Over elevated surfaces, heating of the ground exceeds the heating of the surrounding air at the same altitude above sea level, creating an associated thermal low over the terrain and enhancing any thermal lows that would have otherwise existed,[48][49] and changing the wind circulation of the region."
$p.Slides.Item(25).Shapes.Item(1).TextFrame.TextRange.Text = "Average wind speeds[edit]"
$p.Slides.Item(25).Shapes.Item(2).TextFrame.TextRange.Text = "Average wind speeds[edit]p As described earlier, prevailing and local winds are not spread evenly across the earth, which means that wind speeds also differ by region."
$p.Slides.Item(26).Shapes.Item(1).TextFrame.TextRange.Text = "Wind power density[edit]"
$p.Slides.Item(26).Shapes.Item(2).TextFrame.TextRange.Text = "Wind power density[edit]p Nowadays, a yardstick used to determine the best locations for wind energy development is referred to as wind power density (WPD)."
$p.Slides.Item(27).Shapes.Item(1).TextFrame.TextRange.Text = "Shear[edit]"
$p.Slides.Item(27).Shapes.Item(2).TextFrame.TextRange.Text = "Shear[edit]p Wind shear, sometimes referred to as wind gradient, is a difference in wind speed and direction over a relatively short distance in the Earth's atmosphere."
$p.Slides.Item(28).Shapes.Item(1).TextFrame.TextRange.Text = "Usage[edit]"
$p.Slides.Item(28).Shapes.Item(2).TextFrame.TextRange.Text = "Usage[edit]h3 History[edit]p As a natural force, the wind was often personified as one or more wind gods or as an expression of the supernatural in many cultures."
$p.Slides.Item(29).Shapes.Item(1).TextFrame.TextRange.Text = "History[edit]"
$p.Slides.Item(29).Shapes.Item(2).TextFrame.TextRange.Text = "History[edit]p As a natural force, the wind was often personified as one or more wind gods or as an expression of the supernatural in many cultures."
$p.Slides.Item(30).Shapes.Item(1).TextFrame.TextRange.Text = "Transportation[edit]"
$p.Slides.Item(30).Shapes.Item(2).TextFrame.TextRange.Text = "Transportation[edit]p There are many different forms of sailing ships, but they all have certain basic things in common."
$p.Slides.Item(31).Shapes.Item(1).TextFrame.TextRange.Text = "Power source[edit]"
$p.Slides.Item(31).Shapes.Item(2).TextFrame.TextRange.Text = "Power source[edit]p Historically, the ancient Sinhalese of Anuradhapura and in other cities around Sri Lanka used the monsoon winds to power furnaces as early as 300 BCE."
$p.Slides.Item(32).Shapes.Item(1).TextFrame.TextRange.Text = "Recreation[edit]"
$p.Slides.Item(32).Shapes.Item(2).TextFrame.TextRange.Text = "Recreation[edit]p Wind figures prominently in several popular sports, including recreational hang gliding, hot air ballooning, kite flying, snowkiting, kite landboarding, kite surfing, paragliding, sailing, and windsurfing."
$p.Slides.Item(33).Shapes.Item(1).TextFrame.TextRange.Text = "Role in the natural world[edit]"
$p.Slides.Item(33).Shapes.Item(2).TextFrame.TextRange.Text = "Role in the natural world[edit]p In arid climates, the main source of erosion is wind."
$p.Slides.Item(34).Shapes.Item(1).TextFrame.TextRange.Text = "Erosion[edit]"
$p.Slides.Item(34).Shapes.Item(2).TextFrame.TextRange.Text = "Erosion[edit]p Erosion can be the result of material movement by the wind."
$p.Slides.Item(35).Shapes.Item(1).TextFrame.TextRange.Text = "Desert dust migration[edit]"
$p.Slides.Item(35).Shapes.Item(2).TextFrame.TextRange.Text = "Desert dust migration[edit]p During mid-summer (July in the northern hemisphere), the westward-moving trade winds south of the northward-moving subtropical ridge expand northwestward from the Caribbean into southeastern North America."
$p.Slides.Item(36).Shapes.Item(1).TextFrame.TextRange.Text = "Effect on plants[edit]"
$p.Slides.Item(36).Shapes.Item(2).TextFrame.TextRange.Text = "Effect on plants[edit]p Wind dispersal of seeds, or anemochory, is one of the more primitive means of dispersal."
$p.Slides.Item(37).Shapes.Item(1).TextFrame.TextRange.Text = "Effect on animals[edit]"
$p.Slides.Item(37).Shapes.Item(2).TextFrame.TextRange.Text = "Effect on animals[edit]p Cattle and sheep are prone to wind chill caused by a combination of wind and cold temperatures, when winds exceed 40 kilometers per hour (25 mph), rendering their hair and wool coverings ineffective."
$p.Slides.Item(38).Shapes.Item(1).TextFrame.TextRange.Text = "Sound generation[edit]"
$p.Slides.Item(38).Shapes.Item(2).TextFrame.TextRange.Text = "Sound generation[edit]p Wind causes the generation of sound."
$p.Slides.Item(39).Shapes.Item(1).TextFrame.TextRange.Text = "Related damage[edit]"
$p.Slides.Item(39).Shapes.Item(2).TextFrame.TextRange.Text = "Related damage[edit]p High winds are known to cause damage, depending upon the magnitude of their velocity and pressure differential."

# --- New slides 40-59 (appended at the end, Title and Content layout) ---
$s40 = $p.Slides.Add(40, 2)
$s40.Shapes.Item(1).TextFrame.TextRange.Text = "In outer space[edit]"
$s40.Shapes.Item(2).TextFrame.TextRange.Text = "In outer space[edit]p The solar wind is quite different from a terrestrial wind, in that its origin is the sun, and it is composed of charged particles that have escaped the sun's atmosphere."
$s41 = $p.Slides.Add(41, 2)
$s41.Shapes.Item(1).TextFrame.TextRange.Text = "Planetary wind[edit]"
$s41.Shapes.Item(2).TextFrame.TextRange.Text = "Planetary wind[edit]p The hydrodynamic wind within the upper portion of a planet's atmosphere allows light chemical elements such as hydrogen to move up to the exobase, the lower limit of the exosphere, where the gases can then reach escape velocity, entering outer space without impacting other particles of gas."
$s42 = $p.Slides.Add(42, 2)
$s42.Shapes.Item(1).TextFrame.TextRange.Text = "Solar wind[edit]"
$s42.Shapes.Item(2).TextFrame.TextRange.Text = "Solar wind[edit]p Rather than air, the solar wind is a stream of charged particles—a plasma—ejected from the upper atmosphere of the sun at a rate of 400 kilometers per second (890,000 mph)."
$s43 = $p.Slides.Add(43, 2)
$s43.Shapes.Item(1).TextFrame.TextRange.Text = "On other planets[edit]"
$s43.Shapes.Item(2).TextFrame.TextRange.Text = "On other planets[edit]p Strong 300 kilometers per hour (190 mph) winds at Venus's cloud tops circle the planet every four to five earth days."
$s44 = $p.Slides.Add(44, 2)
$s44.Shapes.Item(1).TextFrame.TextRange.Text = "See also[edit]"
$s45 = $p.Slides.Add(45, 2)
$s45.Shapes.Item(1).TextFrame.TextRange.Text = "References[edit]"
$s46 = $p.Slides.Add(46, 2)
$s46.Shapes.Item(1).TextFrame.TextRange.Text = "External links[edit]"
$s47 = $p.Slides.Add(47, 2)
$s47.Shapes.Item(1).TextFrame.TextRange.Text = "Navigation menu"
$s48 = $p.Slides.Add(48, 2)
$s48.Shapes.Item(1).TextFrame.TextRange.Text = "Personal tools"
$s49 = $p.Slides.Add(49, 2)
$s49.Shapes.Item(1).TextFrame.TextRange.Text = "Namespaces"
$s50 = $p.Slides.Add(50, 2)
$s50.Shapes.Item(1).TextFrame.TextRange.Text = "Variants"
$s51 = $p.Slides.Add(51, 2)
$s51.Shapes.Item(1).TextFrame.TextRange.Text = "Views"
$s52 = $p.Slides.Add(52, 2)
$s52.Shapes.Item(1).TextFrame.TextRange.Text = "More"
$s53 = $p.Slides.Add(53, 2)
$s53.Shapes.Item(1).TextFrame.TextRange.Text = "Search"
$s54 = $p.Slides.Add(54, 2)
$s54.Shapes.Item(1).TextFrame.TextRange.Text = "Navigation"
$s55 = $p.Slides.Add(55, 2)
$s55.Shapes.Item(1).TextFrame.TextRange.Text = "Contribute"
$s56 = $p.Slides.Add(56, 2)
$s56.Shapes.Item(1).TextFrame.TextRange.Text = "Tools"
$s57 = $p.Slides.Add(57, 2)
$s57.Shapes.Item(1).TextFrame.TextRange.Text = "Print/export"
$s58 = $p.Slides.Add(58, 2)
$s58.Shapes.Item(1).TextFrame.TextRange.Text = "In other projects"
$s59 = $p.Slides.Add(59, 2)
$s59.Shapes.Item(1).TextFrame.TextRange.Text = "Languages"

